# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K"), rows 2-33
$newK = @{
    2  = 2
    3  = 3
    4  = 3
    5  = 4
    6  = 1
    7  = 3
    8  = 2
    9  = 5
    10 = 1
    11 = 7
    12 = 4
    13 = 4
    14 = 8
    15 = 2
    16 = 5
    17 = 7
    18 = 3
    19 = 6
    20 = 4
    21 = 7
    22 = 6
    23 = 6
    24 = 3
    25 = 9
    26 = 7
    27 = 6
    28 = 10
    29 = 6
    30 = 6
    31 = 5
    32 = 2
    33 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
